# Applies the cryptos-list price/volume refresh described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, column, new value, and whether the text needs to be forced
# to stay text (so Excel does not silently coerce numeric-looking strings like
# "622.25" or "1.00" into numbers and drop significant trailing zeros).
$updates = @(
    @(2, 4, "69.277.97", $false),
    @(2, 5, "  +2.05%  ", $false),
    @(3, 4, "3.775.32", $false),
    @(3, 5, "  -0.01%  ", $false),
    @(4, 5, "  +0.17%  ", $false),
    @(5, 4, "622.25", $true),
    @(5, 5, "  +4.02%  ", $false),
    @(6, 4, "165.28", $true),
    @(6, 5, "  +1.38%  ", $false),
    @(7, 4, "3.773.88", $false),
    @(7, 5, "  +0.03%  ", $false),
    @(8, 5, "  -0.20%  ", $false),
    @(9, 5, "  +1.50%  ", $false),
    @(10, 4, "0.160", $true),
    @(10, 5, "  +2.47%  ", $false),
    @(11, 5, "  +1.38%  ", $false),
    @(12, 4, "6.65", $true),
    @(12, 5, "  +1.60%  ", $false),
    @(13, 5, "  +0.49%  ", $false),
    @(14, 4, "35.67", $true),
    @(14, 5, "  +1.10%  ", $false),
    @(15, 4, "4.406.07", $false),
    @(15, 5, "  -0.11%  ", $false),
    @(16, 4, "3.709.57", $false),
    @(16, 5, "  -1.75%  ", $false),
    @(17, 4, "69.192.93", $false),
    @(17, 5, "  +1.93%  ", $false),
    @(18, 4, "17.65", $true),
    @(18, 5, "  -3.12%  ", $false),
    @(19, 4, "7.09", $true),
    @(19, 5, "  +1.39%  ", $false),
    @(20, 5, "  -1.07%  ", $false),
    @(21, 4, "468.00", $true),
    @(21, 5, "  +2.12%  ", $false),
    @(22, 4, "9.61", $true),
    @(22, 5, "  +0.22%  ", $false),
    @(23, 4, "0.702", $true),
    @(23, 5, "  +1.17%  ", $false),
    @(24, 5, "  +4.77%  ", $false),
    @(25, 4, "83.30", $true),
    @(25, 5, "  +0.78%  ", $false),
    @(26, 4, "12.00", $true),
    @(26, 5, "  +0.80%  ", $false),
    @(27, 5, "  +3.81%  ", $false),
    @(28, 2, "RenderToken", $false),
    @(28, 3, "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr", $false),
    @(28, 4, "10.01", $true),
    @(28, 5, "  +1.43%  ", $false),
    @(29, 2, "Dai", $false),
    @(29, 3, "https://coinranking.com/coin/MoTuySvg7+dai-dai", $false),
    @(29, 4, "1.00", $true),
    @(29, 5, "  +0.00%  ", $false),
    @(30, 4, "3.921.32", $false),
    @(30, 5, "  -0.16%  ", $false),
    @(31, 2, "PancakeSwap", $false),
    @(31, 3, "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake", $false),
    @(31, 4, "2.66", $true),
    @(31, 5, "  +4.26%  ", $false),
    @(32, 2, "ImmutableX", $false),
    @(32, 3, "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx", $false),
    @(32, 4, "2.24", $true),
    @(32, 5, "  +1.47%  ", $false),
    @(33, 4, "7.30", $true),
    @(33, 5, "  +1.20%  ", $false),
    @(34, 4, "28.81", $true),
    @(35, 4, "1.00", $true),
    @(35, 5, "  -0.06%  ", $false),
    @(36, 4, "3.723.64", $false),
    @(36, 5, "  -0.03%  ", $false),
    @(37, 5, "  +0.68%  ", $false),
    @(38, 5, "  +12.92%  ", $false),
    @(39, 5, "  +3.11%  ", $false),
    @(40, 4, "3.38", $true),
    @(40, 5, "  +6.68%  ", $false),
    @(41, 4, "5.81", $true),
    @(41, 5, "  +0.56%  ", $false),
    @(42, 4, "0.968", $true),
    @(42, 5, "  -1.18%  ", $false),
    @(43, 4, "0.999", $true),
    @(43, 5, "  -0.12%  ", $false),
    @(45, 5, "  +1.74%  ", $false),
    @(46, 4, "154.33", $true),
    @(46, 5, "  +1.07%  ", $false),
    @(47, 4, "43.17", $true),
    @(47, 5, "  -0.41%  ", $false),
    @(48, 5, "  -0.96%  ", $false),
    @(49, 5, "  +3.44%  ", $false),
    @(50, 4, "8.40", $true),
    @(50, 5, "  +1.40%  ", $false),
    @(51, 5, "  +0.03%  ", $false),
)

foreach ($u in $updates) {
    $cell = $ws.Cells.Item($u[0], $u[1])
    if ($u[3]) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $u[2]
}
